$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2748.7
$ws.Range("I28").Value = 947.625
$ws.Range("K28").Value = 947.625
$ws.Range("M28").Value = -462.625
$ws.Range("H64").Value = 90915780
$ws.Range("J64").Value = 166673330
$ws.Range("L64").Value = 166673330
$ws.Range("N64").Value = -166673826
$ws.Range("H67").Value = 90915780
$ws.Range("J67").Value = 166673330
$ws.Range("L67").Value = 166673330
$ws.Range("N67").Value = -166675046
$ws.Range("H107").Value = 616.4
$ws.Range("I107").Value = 623.5789
$ws.Range("K107").Value = 623.5789
$ws.Range("M107").Value = 1296.4211
$ws.Range("H138").Value = 2927.9443
$ws.Range("J138").Value = 2893.8809
$ws.Range("L138").Value = 8681.6427
$ws.Range("N138").Value = -18961.6427
$ws.Range("H141").Value = 5565.263
$ws.Range("I141").Value = 2783.875
$ws.Range("K141").Value = 8351.625
$ws.Range("M141").Value = -3171.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3152.0847
$ws.Range("I32").Value = 2947.8276
$ws.Range("J32").Value = 14999
$ws.Range("K32").Value = 2947.8276
$ws.Range("L32").Value = 14999
$ws.Range("M32").Value = -2660.8276
$ws.Range("N32").Value = -15573
$ws.Range("H45").Value = 86882.8
$ws.Range("I45").Value = 86882.8
$ws.Range("K45").Value = 86882.8
$ws.Range("M45").Value = -86505.8
$ws.Range("H61").Value = 2531.0212
$ws.Range("I61").Value = 2071.3333
$ws.Range("J61").Value = 3614.5715
$ws.Range("K61").Value = 2071.3333
$ws.Range("L61").Value = 3614.5715
$ws.Range("M61").Value = -1859.3333
$ws.Range("N61").Value = -4038.5715
$ws.Range("H122").Value = 4290.5
$ws.Range("I122").Value = 4217.346
$ws.Range("K122").Value = 12652.038
$ws.Range("M122").Value = -10202.038
$ws.Range("H132").Value = 2615
$ws.Range("I132").Value = 1595.3889
$ws.Range("J132").Value = 4654.222
$ws.Range("K132").Value = 4786.1667
$ws.Range("L132").Value = 13962.666
$ws.Range("M132").Value = -2256.1667
$ws.Range("N132").Value = -19022.666
$ws.Range("H136").Value = 2531.0212
$ws.Range("I136").Value = 2071.3333
$ws.Range("J136").Value = 3614.5715
$ws.Range("K136").Value = 6213.999899999999
$ws.Range("L136").Value = 10843.7145
$ws.Range("M136").Value = -3663.999899999999
$ws.Range("N136").Value = -15943.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 329.91666
$ws.Range("I80").Value = 113.25
$ws.Range("J80").Value = 438.25
$ws.Range("K80").Value = 113.25
$ws.Range("L80").Value = 438.25
$ws.Range("M80").Value = 884.75
$ws.Range("N80").Value = -2434.25
$ws.Range("H83").Value = 329.91666
$ws.Range("I83").Value = 113.25
$ws.Range("J83").Value = 438.25
$ws.Range("K83").Value = 566.25
$ws.Range("L83").Value = 2191.25
$ws.Range("M83").Value = 4425.75
$ws.Range("N83").Value = -12175.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1774.8572
$ws.Range("I16").Value = 1839.2142
$ws.Range("J16").Value = 1646.1428
$ws.Range("K16").Value = 1839.2142
$ws.Range("L16").Value = 1646.1428
$ws.Range("M16").Value = -1552.2142
$ws.Range("N16").Value = -2220.1428
$ws.Range("H31").Value = 3931.1836
$ws.Range("I31").Value = 2927.8518
$ws.Range("J31").Value = 5162.5454
$ws.Range("K31").Value = 2927.8518
$ws.Range("L31").Value = 5162.5454
$ws.Range("M31").Value = -2632.8518
$ws.Range("N31").Value = -5752.5454
$ws.Range("H34").Value = 3931.1836
$ws.Range("I34").Value = 2927.8518
$ws.Range("J34").Value = 5162.5454
$ws.Range("K34").Value = 2927.8518
$ws.Range("L34").Value = 5162.5454
$ws.Range("M34").Value = -2725.8518
$ws.Range("N34").Value = -5566.5454
$ws.Range("H58").Value = 2960.3667
$ws.Range("I58").Value = 2522.4167
$ws.Range("J58").Value = 3252.3333
$ws.Range("K58").Value = 2522.4167
$ws.Range("L58").Value = 3252.3333
$ws.Range("M58").Value = -2319.4167
$ws.Range("N58").Value = -3658.3333
$ws.Range("H86").Value = 4998.5
$ws.Range("J86").Value = 4999
$ws.Range("L86").Value = 4999
$ws.Range("N86").Value = -7245
$ws.Range("H89").Value = 4998.5
$ws.Range("J89").Value = 4999
$ws.Range("L89").Value = 24995
$ws.Range("N89").Value = -36227
$ws.Range("H113").Value = 1774.8572
$ws.Range("I113").Value = 1839.2142
$ws.Range("J113").Value = 1646.1428
$ws.Range("K113").Value = 1839.2142
$ws.Range("L113").Value = 1646.1428
$ws.Range("M113").Value = 330.7858000000001
$ws.Range("N113").Value = -5986.1428
$ws.Range("H136").Value = 2960.3667
$ws.Range("I136").Value = 2522.4167
$ws.Range("J136").Value = 3252.3333
$ws.Range("K136").Value = 7567.250100000001
$ws.Range("L136").Value = 9756.999899999999
$ws.Range("M136").Value = -5017.250100000001
$ws.Range("N136").Value = -14856.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2584.4614
$ws.Range("I34").Value = 1649.5
$ws.Range("J34").Value = 2754.4546
$ws.Range("K34").Value = 4948.5
$ws.Range("L34").Value = 8263.363799999999
$ws.Range("M34").Value = -4864.5
$ws.Range("N34").Value = -8431.363799999999
$ws.Range("H101").Value = 7833.3335
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 7833.3335
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 23500.0005
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -28368.0005
$ws.Range("H104").Value = 6000
$ws.Range("J104").Value = 6000
$ws.Range("L104").Value = 18000
$ws.Range("N104").Value = -23242
$ws.Range("H132").Value = 2293.6365
$ws.Range("J132").Value = 3489.5
$ws.Range("L132").Value = 31405.5
$ws.Range("N132").Value = -36465.5
$ws.Range("H141").Value = 12066.357
$ws.Range("I141").Value = 4448.091
$ws.Range("K141").Value = 13344.273
$ws.Range("M141").Value = -8164.273000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 8046.75
$ws.Range("I99").Value = 2222.6667
$ws.Range("J99").Value = 25519
$ws.Range("K99").Value = 2222.6667
$ws.Range("L99").Value = 25519
$ws.Range("M99").Value = 23.33329999999978
$ws.Range("N99").Value = -30011
$ws.Range("H112").Value = 92822
$ws.Range("J112").Value = 92822
$ws.Range("L112").Value = 92822
$ws.Range("N112").Value = -95038
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1106.4286
$ws.Range("J22").Value = 1299
$ws.Range("L22").Value = 1299
$ws.Range("N22").Value = -1889
$ws.Range("H27").Value = 1106.4286
$ws.Range("J27").Value = 1299
$ws.Range("L27").Value = 1299
$ws.Range("N27").Value = -1513
$ws.Range("H56").Value = 31008.2
$ws.Range("I56").Value = 20347
$ws.Range("K56").Value = 20347
$ws.Range("M56").Value = -19656
$ws.Range("H61").Value = 2777.9412
$ws.Range("I61").Value = 2751.6428
$ws.Range("K61").Value = 2751.6428
$ws.Range("M61").Value = -2549.6428
$ws.Range("H99").Value = 24991
$ws.Range("I99").Value = 24991
$ws.Range("K99").Value = 24991
$ws.Range("M99").Value = -21996
$ws.Range("H113").Value = 2777.9412
$ws.Range("I113").Value = 2751.6428
$ws.Range("K113").Value = 2751.6428
$ws.Range("M113").Value = -581.6428000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 15500
$ws.Range("J51").Value = 35000
$ws.Range("L51").Value = 35000
$ws.Range("N51").Value = -36020
$ws.Range("H132").Value = 3897.9473
$ws.Range("I132").Value = 3791.4
$ws.Range("J132").Value = 4297.5
$ws.Range("K132").Value = 11374.2
$ws.Range("L132").Value = 12892.5
$ws.Range("M132").Value = -8844.200000000001
$ws.Range("N132").Value = -17952.5
$ws.Range("H136").Value = 55559810
$ws.Range("I136").Value = 71429820
$ws.Range("J136").Value = 14748.5
$ws.Range("K136").Value = 214289460
$ws.Range("L136").Value = 44245.5
$ws.Range("M136").Value = -214286910
$ws.Range("N136").Value = -49345.5
